# Applies the textual edits described by the commit
# "updatera ratta fel pa olika klasser" (update/fix errors on various classes)
# to the Redovisning rapport document.

$d = $word.ActiveDocument

# 1. Cam's paragraph: add "testa applikationen" after "...kodning av GUI,"
$d.Content.Find.Execute(
    "databasscript, delaktig i design och kodning av GUI, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "databasscript, delaktig i design och kodning av GUI, testa applikationen ", 2)

# 2. Hampus's paragraph: "testning av applikationen" -> "test av applikationen"
$d.Content.Find.Execute(
    "testning av applikationen",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "test av applikationen", 2)

# 3. Asa's paragraph: ")???" -> ")???, testa applikationen"
$d.Content.Find.Execute(
    ")???",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ")???, testa applikationen", 2)

# 4. Arbetsformer paragraph: "individ" -> "medlem"
$d.Content.Find.Execute(
    ", varje individ få",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", varje medlem få", 2)

# 5. "programmeringsuppgifter" -> "programmerings- uppgifter"
$d.Content.Find.Execute(
    "programmeringsuppgifter",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "programmerings- uppgifter", 2)

# 6. Append new sentence after "pa egen hand."
$d.Content.Find.Execute(
    "på egen hand.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "på egen hand. Sedan få varje medlem bestämma själv om dem vill samarbeta med andra.", 2)

# 7. "pdf dokument. " -> "pdf dokumenten. "
$d.Content.Find.Execute(
    "pdf dokument. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "pdf dokumenten. ", 2)
